$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '55.697.53'
$ws.Range("E2").Value = '  +8.01%  '

$ws.Range("D3").Value = '3.211.53'
$ws.Range("E3").Value = '  +3.37%  '

$ws.Range("E4").Value = '  -0.01%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '396.45'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +2.56%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '109.57'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +6.18%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.554'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +2.68%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.999'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -0.10%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.621'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +6.04%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '39.23'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +5.87%  '

$ws.Range("E11").Value = '  +5.27%  '

$ws.Range("E12").Value = '  +2.15%  '

$ws.Range("D13").Value = '3.717.11'
$ws.Range("E13").Value = '  +3.48%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '19.02'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +2.35%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '8.05'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +2.65%  '

$ws.Range("B16").Value = 'Polygon'
$ws.Range("C16").Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '1.05'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +5.79%  '

$ws.Range("B17").Value = 'WrappedEther'
$ws.Range("C17").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D17").Value = '3.190.94'
$ws.Range("E17").Value = '  +2.65%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '10.59'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -3.37%  '

$ws.Range("D19").Value = '55.604.69'
$ws.Range("E19").Value = '  +7.77%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '3.38'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +3.55%  '

$ws.Range("E21").Value = '  +5.57%  '

$ws.Range("E22").Value = '  +5.55%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '306.35'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +14.79%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '75.42'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +7.87%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '3.22'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +1.47%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '8.27'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +1.82%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '28.16'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +4.18%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '7.52'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +4.02%  '

$ws.Range("E29").Value = '  +4.34%  '

$ws.Range("E30").Value = '  -0.11%  '

$ws.Range("B31").Value = 'Hedera'
$ws.Range("C31").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.110'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +3.74%  '

$ws.Range("B32").Value = 'Cosmos'
$ws.Range("C32").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '11.30'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +8.94%  '

$ws.Range("E33").Value = '  +2.71%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '36.11'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +2.52%  '

$ws.Range("E35").Value = '  +2.32%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '51.20'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +2.19%  '

$ws.Range("B37").Value = 'FirstDigitalUSD'
$ws.Range("C37").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.00'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +0.19%  '

$ws.Range("B38").Value = 'Stacks'
$ws.Range("C38").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '3.09'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +23.41%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '3.49'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +4.16%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '133.59'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +3.33%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '4.04'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +8.42%  '

$ws.Range("E42").Value = '  +2.09%  '

$ws.Range("B43").Value = 'TheGraph'
$ws.Range("C43").Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.286'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -2.13%  '

$ws.Range("B44").Value = 'Stellar'
$ws.Range("C44").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.119'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +2.80%  '

$ws.Range("B45").Value = 'Celestia'
$ws.Range("C45").Value = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '17.02'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +2.48%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '22.19'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -1.07%  '

$ws.Range("B47").Value = 'ThetaToken'
$ws.Range("C47").Value = 'https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '2.13'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +45.88%  '

$ws.Range("B48").Value = 'WEMIXToken'
$ws.Range("C48").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '2.11'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +1.56%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '2.48'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -1.17%  '

$ws.Range("B50").Value = 'Maker'
$ws.Range("C50").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D50").Value = '2.141.36'
$ws.Range("E50").Value = '  +3.58%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.0363'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +8.83%  '

